$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.089.17"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.661.65"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.39%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "208.24"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.5183"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.29%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.2584"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -1.82%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.06294"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("E10").Value = "  +0.61%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.07539"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.660.09"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -0.33%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.5388"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -3.15%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "66.13"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "0.0₅7927"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "26.096.24"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  -0.26%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "4.699"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -1.82%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "187.41"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("E21").Value = "  -1.85%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "6.195"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  -0.45%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "148.06"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.1212"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -2.45%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "7.379"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "15.64"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "1.388"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "0.06014"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -4.48%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.263"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -0.16%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.399"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -0.47%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.633"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.9846"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "2.759"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "2.387"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "1.104.81"
$ws.Range("E38").Value = "  +1.14%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.01591"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "5.968"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -2.21%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.8479"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  -0.11%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "99.92"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "1.815.20"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "0.0₈107"
$ws.Range("E45").Value = "  +0.46%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "55.08"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.63%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "8.018"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.05228"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.4239"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "5.854"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
